# Weekly update: insert a new (most recent) price record at row 140 for
# "Agrícola del Norte S.A. de Arica - Locoto", pushing the existing
# rows 140:187 down to 141:188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 140 (shifts rows 140-187 down to 141-188,
# carrying their formatting/styles with them).
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with this week's record. The non-varying
# columns (market/region/category/variety/unit/origin/kg) are constant
# across the whole sheet.
$row = 140
$ws.Cells.Item($row, 1).Value2  = 1
$ws.Cells.Item($row, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value2  = 45229
$ws.Cells.Item($row, 5).Value2  = 15
$ws.Cells.Item($row, 6).Value2  = 100112042
$ws.Cells.Item($row, 7).Value2  = "Locoto"
$ws.Cells.Item($row, 8).Value2  = "Sin especificar"
$ws.Cells.Item($row, 9).Value2  = "Primera"
$ws.Cells.Item($row, 10).Value2 = 130
$ws.Cells.Item($row, 11).Value2 = 48000
$ws.Cells.Item($row, 12).Value2 = 50000
$ws.Cells.Item($row, 13).Value2 = 48923
$ws.Cells.Item($row, 14).Value2 = "$/caja 20 kilos"
$ws.Cells.Item($row, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value2 = 2446
$ws.Cells.Item($row, 17).Value2 = 20
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
